$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price table rows with the latest scraped values.
# Column D ("Price") values are forced to text with a leading apostrophe
# so Excel does not reinterpret dotted/decimal-looking strings as numbers.

$ws.Range("D2").Value = "'26.663.41"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "'1.598.69"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'211.62"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "'0.247"
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").Value = "'19.56"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "'1.822.43"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "'1.597.97"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "'65.12"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "'26.638.63"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'209.26"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.00"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'7.04"
$ws.Range("E21").Value = "  +4.10%  "
$ws.Range("D22").Value = "'4.29"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").Value = "'8.99"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").Value = "'144.44"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "'7.11"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").Value = "'15.29"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").Value = "'3.25"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "'2.96"
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("D34").Value = "'1.281.61"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").Value = "'0.619"
$ws.Range("E35").Value = "  -7.50%  "
$ws.Range("D36").Value = "'2.45"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").Value = "'0.834"
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("E40").Value = "  +17.09%  "
$ws.Range("D41").Value = "'5.48"
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "'0.784"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").Value = "'63.45"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "'1.733.40"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "'1.57"
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  -1.76%  "
